$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.742.05"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "3.148.84"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.90"
$ws.Range("E5").Value = "  +1.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.10"
$ws.Range("E6").Value = "  -0.91%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "3.149.74"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  -0.33%  "

$ws.Range("E10").Value = "  -2.29%  "

$ws.Range("E11").Value = "  -0.77%  "

$ws.Range("E12").Value = "  -0.68%  "

$ws.Range("E13").Value = "  +2.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.12"
$ws.Range("E14").Value = "  -2.64%  "

$ws.Range("D15").Value = "3.666.69"
$ws.Range("E15").Value = "  -0.16%  "

$ws.Range("D16").Value = "64.809.64"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.159.89"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.15"
$ws.Range("E18").Value = "  -1.17%  "

$ws.Range("E19").Value = "  +0.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "504.87"
$ws.Range("E20").Value = "  -2.31%  "

$ws.Range("E21").Value = "  -0.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.714"
$ws.Range("E22").Value = "  -3.04%  "

$ws.Range("E23").Value = "  -0.75%  "

$ws.Range("E24").Value = "  -1.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.43"
$ws.Range("E25").Value = "  -0.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("E27").Value = "  +2.20%  "

$ws.Range("E28").Value = "  -0.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.78"
$ws.Range("E30").Value = "  +3.76%  "

$ws.Range("E31").Value = "  -1.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("E33").Value = "  +0.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.25"
$ws.Range("E34").Value = "  +2.15%  "

$ws.Range("E35").Value = "  -1.42%  "

$ws.Range("E36").Value = "  -1.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "483.11"
$ws.Range("E37").Value = "  -1.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0892"
$ws.Range("E38").Value = "  +2.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0415"
$ws.Range("E39").Value = "  -1.87%  "

$ws.Range("E40").Value = "  -1.87%  "

$ws.Range("E41").Value = "  +0.91%  "

$ws.Range("D42").Value = "2.987.58"
$ws.Range("E42").Value = "  -4.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.118"
$ws.Range("E43").Value = "  -1.97%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.42"
$ws.Range("E44").Value = "  -0.86%  "

$ws.Range("E45").Value = "  -4.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.19"
$ws.Range("E46").Value = "  -3.82%  "

$ws.Range("D47").Value = "0.0₃0587"
$ws.Range("E47").Value = "  +1.92%  "

$ws.Range("E49").Value = "  -1.53%  "

$ws.Range("E50").Value = "  -2.77%  "

$ws.Range("E51").Value = "  +14.29%  "
